$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 16. This shifts the existing
# rows 16-36 down to become rows 18-38, matching the new dimension A1:T38.
$ws.Rows("16:17").Insert()

# Fill in the two newly inserted rows (16 and 17) with the new weekly data.

# Row 16
$ws.Cells.Item(16, 1).Value = 6
$ws.Cells.Item(16, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(16, 3).Value = "Metropolitana"
$ws.Cells.Item(16, 4).Value = 44672
$ws.Cells.Item(16, 5).Value = 13
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100102
$ws.Cells.Item(16, 8).Value = "Cítricos"
$ws.Cells.Item(16, 9).Value = 100102006
$ws.Cells.Item(16, 10).Value = "Pomelo"
$ws.Cells.Item(16, 11).Value = "Start Ruby"
$ws.Cells.Item(16, 12).Value = "Especial"
$ws.Cells.Item(16, 13).Value = 8
$ws.Cells.Item(16, 14).Value = 200000
$ws.Cells.Item(16, 15).Value = 200000
$ws.Cells.Item(16, 16).Value = 200000
$ws.Cells.Item(16, 17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(16, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(16, 19).Value = 571
$ws.Cells.Item(16, 20).Value = 350

# Row 17
$ws.Cells.Item(17, 1).Value = 6
$ws.Cells.Item(17, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(17, 3).Value = "Metropolitana"
$ws.Cells.Item(17, 4).Value = 44672
$ws.Cells.Item(17, 5).Value = 13
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100102
$ws.Cells.Item(17, 8).Value = "Cítricos"
$ws.Cells.Item(17, 9).Value = 100102006
$ws.Cells.Item(17, 10).Value = "Pomelo"
$ws.Cells.Item(17, 11).Value = "Start Ruby"
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 10
$ws.Cells.Item(17, 14).Value = 150000
$ws.Cells.Item(17, 15).Value = 150000
$ws.Cells.Item(17, 16).Value = 150000
$ws.Cells.Item(17, 17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(17, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(17, 19).Value = 429
$ws.Cells.Item(17, 20).Value = 350
